$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (Comments / MU) -------------------------------------
$ws.Cells.Item(1, 15).Value = "Comments"
$ws.Cells.Item(1, 16).Value = "MU"

# --- Row 2 edits -------------------------------------------------------
$ws.Cells.Item(2, 1).Value = 263
$ws.Cells.Item(2, 2).Value = "NL39871W4D2"
$ws.Cells.Item(2, 3).Value = "14WN1699"
$ws.Cells.Item(2, 5).Value = "14K WG BAGUETTE DIA TASSLE NECKLACE 16""+2"" CHAIN"
$ws.Cells.Item(2, 6).ClearContents()
$ws.Cells.Item(2, 7).Value = 6.006
$ws.Cells.Item(2, 8).Value = 0.78
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 5.85
$ws.Cells.Item(2, 11).Value = 69
$ws.Cells.Item(2, 12).Value = 3290
$ws.Cells.Item(2, 13).Value = 1
$ws.Cells.Item(2, 14).Value = 589

# --- Row 3 (new) ---------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 82.5
$ws.Cells.Item(3, 1).Value = 258
$ws.Cells.Item(3, 2).Value = "OCKSFORHVX530/4"
$ws.Cells.Item(3, 3).Value = "14YC361"
$ws.Cells.Item(3, 5).Value = "14K YG Swivel Lobster fancy clasp 17.5x7.7mm"
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 1.02
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 50
$ws.Cells.Item(3, 13).Value = 51
$ws.Cells.Item(3, 14).Value = 5810
$ws.Cells.Item(3, 15).Value = "Italian clasp"

# --- Row 4 (new) ---------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 82.5
$ws.Cells.Item(4, 1).Value = 101
$ws.Cells.Item(4, 2).Value = "V101"
$ws.Cells.Item(4, 3).Value = "SAM1"
$ws.Cells.Item(4, 5).Value = "Hello"
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1
$ws.Cells.Item(4, 14).Value = 581
$ws.Cells.Item(4, 15).Value = "Heya"
$ws.Cells.Item(4, 16).Value = 14

# --- New column width for Comments column (O) -----------------------------
$ws.Columns.Item(15).ColumnWidth = 25
